# Daily attendance processing - 2025-12-08 20:30:36
# Normalize the "Recorded By" (column G) values so that the "System" token
# is moved to the end of the comma-separated list, preserving the relative
# order of the other tokens (and preserving duplicate "System" tokens, if any).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if (-not $val.Contains("System")) { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $nonSystem = @()
    $systemCount = 0
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) {
            $systemCount += 1
        } else {
            $nonSystem += $p
        }
    }

    if ($systemCount -eq 0) { continue }

    $newParts = @()
    $newParts += $nonSystem
    for ($i = 0; $i -lt $systemCount; $i++) { $newParts += "System" }

    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
